$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.134130716323853
$ws.Range("B1").Value = 2.221505880355835
$ws.Range("C1").Value = 10.79817676544189
$ws.Range("D1").Value = 2.288745641708374
$ws.Range("E1").Value = 1.281835317611694
